$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''68.837.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +1.60%  '
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = '''3.864.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.75%  '
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = '''  -0.09%  '
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = '''602.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +1.02%  '
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = '''171.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +3.87%  '
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = '''3.863.82'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +0.75%  '
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = '''  +0.02%  '
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = '''0.530'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +1.24%  '
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = '''  +3.63%  '
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = '''  +4.24%  '
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = '''  +1.80%  '
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = '''0.0000287'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +16.80%  '
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = '''37.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +1.64%  '
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = '''4.512.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +0.72%  '
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = '''3.877.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +1.40%  '
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = '''68.826.48'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +1.62%  '
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = '''18.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +1.10%  '
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = '''7.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +0.27%  '
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = '''  +0.63%  '
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = '''11.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +4.57%  '
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = '''473.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +1.69%  '
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = '''  +0.96%  '
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = '''0.0000163'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +3.38%  '
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = '''83.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +1.16%  '
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = '''2.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +3.50%  '
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = '''12.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +0.87%  '
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = '''10.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +5.20%  '
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = '''  +0.02%  '
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = '''2.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +1.47%  '
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = '''4.017.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.82%  '
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = '''7.80'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +1.73%  '
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("B33").Value = '''ImmutableX'
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = '''2.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +1.50%  '
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("B34").Value = '''EthereumClassic'
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = '''31.38'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +2.02%  '
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = '''  +1.74%  '
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = '''3.828.99'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +0.55%  '
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = '''  +25.17%  '
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = '''  +1.14%  '
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = '''1.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +1.81%  '
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("B40").Value = '''Filecoin'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = '''5.99'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +2.28%  '
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = '''Kaspa'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = '''0.140'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +0.66%  '
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = '''  +0.12%  '
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").Value = '''  +3.23%  '
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").Value = '''FLOKI'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = '''https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = '''0.000304'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +15.84%  '
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("B45").Value = '''Stacks'
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = '''2.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +1.87%  '
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = '''USDe'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = '''https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = '''1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +0.03%  '
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("B47").Value = '''Cosmos'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''8.73'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +2.76%  '
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("B48").Value = '''Bittensor'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = '''419.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -0.55%  '
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = '''46.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -1.56%  '
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = '''142.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -0.36%  '
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = '''0.0359'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +1.62%  '
$ws.Range("E51").Style = "Normal"
